$wb = $excel.ActiveWorkbook

# ---- Sheet1: Overview ----
$ws = $wb.Worksheets.Item(1)
# Plain (non-hyperlinked) cell values
$ws.Range('A1').Value2 = 'File Name'
$ws.Range('B1').Value2 = 'zh-cn'
$ws.Range('C1').Value2 = 'de-de'
$ws.Range('D1').Value2 = 'Latest Handoff Date'
$ws.Range('B2').Value2 = 'Handed back: in sync with en-US'
$ws.Range('C2').Value2 = 'Handed back: in sync with en-US'
$ws.Range('D2').Value2 = '2016-03-24 08:49:21'
$ws.Range('B3').Value2 = 'Handed back: in sync with en-US'
$ws.Range('C3').Value2 = 'Handed back: in sync with en-US'
$ws.Range('D3').Value2 = '2016-03-24 08:49:21'
$ws.Range('B4').Value2 = 'Ready for handoff'
$ws.Range('C4').Value2 = 'Ready for handoff'
$ws.Range('D4').Value2 = '2016-03-24 08:52:46'
# Rebuild hyperlinks (clears stale ones, then re-adds in original rId order
# with the updated display text so cell values + link labels match the refreshed report)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/2976a60940be99f39c269d8b6da1f596b51c912d/e2e/2f4bd7cb-c696-4503-99c3-da8190232af0.md', "", "", 'ffff11cb6873-c225-4731-8417-48de30a3441c.md')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/2976a60940be99f39c269d8b6da1f596b51c912d/e2e/ffff11cb6873-c225-4731-8417-48de30a3441c.md', "", "", 'ffffff9de48753-ae7a-44ad-bb99-49002869434e.md')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/2976a60940be99f39c269d8b6da1f596b51c912d/e2e/ffffff9de48753-ae7a-44ad-bb99-49002869434e.md', "", "", '2f4bd7cb-c696-4503-99c3-da8190232af0.md')

# ---- Sheet2: zh-cn ----
$ws = $wb.Worksheets.Item(2)
# Plain (non-hyperlinked) cell values
$ws.Range('A1').Value2 = 'Source File Name'
$ws.Range('B1').Value2 = 'File Extension'
$ws.Range('C1').Value2 = 'Status'
$ws.Range('D1').Value2 = 'Latest Handoff File'
$ws.Range('E1').Value2 = 'Latest Handoff Datetime'
$ws.Range('F1').Value2 = 'Latest Target File'
$ws.Range('G1').Value2 = 'Latest Handback File'
$ws.Range('H1').Value2 = 'Latest Handback DateTime'
$ws.Range('J1').Value2 = 'Reference Tokens'
$ws.Range('K1').Value2 = 'Handoff Reason'
$ws.Range('L1').Value2 = 'Dependency From'
$ws.Range('B2').Value2 = '.md'
$ws.Range('C2').Value2 = 'Handed back: in sync with en-US'
$ws.Range('E2').Value2 = '2016-03-24 08:49:16'
$ws.Range('H2').Value2 = '2016-03-24 08:49:45'
$ws.Range('J2').Value2 = 'Include'
$ws.Range('B3').Value2 = '.md'
$ws.Range('C3').Value2 = 'Handed back: in sync with en-US'
$ws.Range('E3').Value2 = '2016-03-24 08:49:16'
$ws.Range('H3').Value2 = '2016-03-24 08:49:45'
$ws.Range('J3').Value2 = 'Include'
$ws.Range('B4').Value2 = '.md'
$ws.Range('C4').Value2 = 'Ready for handoff'
$ws.Range('E4').Value2 = '2016-03-24 08:52:42'
$ws.Range('H4').Value2 = '2016-03-24 08:51:49'
$ws.Range('J4').Value2 = 'Include'
# Rebuild hyperlinks (clears stale ones, then re-adds in original rId order
# with the updated display text so cell values + link labels match the refreshed report)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/2976a60940be99f39c269d8b6da1f596b51c912d/e2e/2f4bd7cb-c696-4503-99c3-da8190232af0.md', "", "", 'ffff11cb6873-c225-4731-8417-48de30a3441c.md')
$ws.Hyperlinks.Add($ws.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dcdc7e9e2caef0f8c342471d89a1be698d7c26b8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2f4bd7cb-c696-4503-99c3-da8190232af0.b64a3c588d434ab4c7a8a27f5428149657c67613.zh-cn.xlf', "", "", '34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/8b27d7bc354527765a27f1103b8caa10e1a94c92/e2e/2f4bd7cb-c696-4503-99c3-da8190232af0.md', "", "", '34ded686-4006-40a0-a24a-57ef94237596.md')
$ws.Hyperlinks.Add($ws.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/8e1577156e2f38c9ed359369d04ec717e2b99dc5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2f4bd7cb-c696-4503-99c3-da8190232af0.b64a3c588d434ab4c7a8a27f5428149657c67613.zh-cn.xlf', "", "", '34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/2976a60940be99f39c269d8b6da1f596b51c912d/e2e/ffff11cb6873-c225-4731-8417-48de30a3441c.md', "", "", 'ffffff9de48753-ae7a-44ad-bb99-49002869434e.md')
$ws.Hyperlinks.Add($ws.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/72963e108308c1c39870dac2d05435e3b1434867/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf', "", "", '34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6f60f6f1fa63f4fa6871ea33c017c179aa8f9138/e2e/34ded686-4006-40a0-a24a-57ef94237596.md', "", "", '34ded686-4006-40a0-a24a-57ef94237596.md')
$ws.Hyperlinks.Add($ws.Range('G3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/a095b6d4a6367fee06428c929a72f6d8e2c55e54/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf', "", "", '34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/2976a60940be99f39c269d8b6da1f596b51c912d/e2e/ffffff9de48753-ae7a-44ad-bb99-49002869434e.md', "", "", '2f4bd7cb-c696-4503-99c3-da8190232af0.md')
$ws.Hyperlinks.Add($ws.Range('D4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/72963e108308c1c39870dac2d05435e3b1434867/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf', "", "", '2f4bd7cb-c696-4503-99c3-da8190232af0.b64a3c588d434ab4c7a8a27f5428149657c67613.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6f60f6f1fa63f4fa6871ea33c017c179aa8f9138/e2e/34ded686-4006-40a0-a24a-57ef94237596.md', "", "", '2f4bd7cb-c696-4503-99c3-da8190232af0.md')
$ws.Hyperlinks.Add($ws.Range('G4'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/a095b6d4a6367fee06428c929a72f6d8e2c55e54/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.zh-cn.xlf', "", "", '2f4bd7cb-c696-4503-99c3-da8190232af0.b64a3c588d434ab4c7a8a27f5428149657c67613.zh-cn.xlf')

# ---- Sheet3: de-de ----
$ws = $wb.Worksheets.Item(3)
# Plain (non-hyperlinked) cell values
$ws.Range('A1').Value2 = 'Source File Name'
$ws.Range('B1').Value2 = 'File Extension'
$ws.Range('C1').Value2 = 'Status'
$ws.Range('D1').Value2 = 'Latest Handoff File'
$ws.Range('E1').Value2 = 'Latest Handoff Datetime'
$ws.Range('F1').Value2 = 'Latest Target File'
$ws.Range('G1').Value2 = 'Latest Handback File'
$ws.Range('H1').Value2 = 'Latest Handback DateTime'
$ws.Range('J1').Value2 = 'Reference Tokens'
$ws.Range('K1').Value2 = 'Handoff Reason'
$ws.Range('L1').Value2 = 'Dependency From'
$ws.Range('B2').Value2 = '.md'
$ws.Range('C2').Value2 = 'Handed back: in sync with en-US'
$ws.Range('E2').Value2 = '2016-03-24 08:49:21'
$ws.Range('H2').Value2 = '2016-03-24 08:49:52'
$ws.Range('J2').Value2 = 'Include'
$ws.Range('B3').Value2 = '.md'
$ws.Range('C3').Value2 = 'Handed back: in sync with en-US'
$ws.Range('E3').Value2 = '2016-03-24 08:49:21'
$ws.Range('H3').Value2 = '2016-03-24 08:49:52'
$ws.Range('J3').Value2 = 'Include'
$ws.Range('B4').Value2 = '.md'
$ws.Range('C4').Value2 = 'Ready for handoff'
$ws.Range('E4').Value2 = '2016-03-24 08:52:46'
$ws.Range('H4').Value2 = '2016-03-24 08:51:57'
$ws.Range('J4').Value2 = 'Include'
# Rebuild hyperlinks (clears stale ones, then re-adds in original rId order
# with the updated display text so cell values + link labels match the refreshed report)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/2976a60940be99f39c269d8b6da1f596b51c912d/e2e/2f4bd7cb-c696-4503-99c3-da8190232af0.md', "", "", 'ffff11cb6873-c225-4731-8417-48de30a3441c.md')
$ws.Hyperlinks.Add($ws.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67a2107f250953ed45b636cd46deb3aa258395de/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2f4bd7cb-c696-4503-99c3-da8190232af0.b64a3c588d434ab4c7a8a27f5428149657c67613.de-de.xlf', "", "", '34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/68f514c6724084feb18dcb2c0af2e9e28241c41a/e2e/2f4bd7cb-c696-4503-99c3-da8190232af0.md', "", "", '34ded686-4006-40a0-a24a-57ef94237596.md')
$ws.Hyperlinks.Add($ws.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/3d82cfc4a27ac6de71bcc51770ed2336faa39890/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2f4bd7cb-c696-4503-99c3-da8190232af0.b64a3c588d434ab4c7a8a27f5428149657c67613.de-de.xlf', "", "", '34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/2976a60940be99f39c269d8b6da1f596b51c912d/e2e/ffff11cb6873-c225-4731-8417-48de30a3441c.md', "", "", 'ffffff9de48753-ae7a-44ad-bb99-49002869434e.md')
$ws.Hyperlinks.Add($ws.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c32d239dcec73e8605290148acaa033669ff389f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf', "", "", '34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/54a8ca9bebda7d35e8d6e69eac6b9ed3df74a6c7/e2e/34ded686-4006-40a0-a24a-57ef94237596.md', "", "", '34ded686-4006-40a0-a24a-57ef94237596.md')
$ws.Hyperlinks.Add($ws.Range('G3'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/8d1000a3019e2fd70ed23b4162a640027ce1539d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf', "", "", '34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/2976a60940be99f39c269d8b6da1f596b51c912d/e2e/ffffff9de48753-ae7a-44ad-bb99-49002869434e.md', "", "", '2f4bd7cb-c696-4503-99c3-da8190232af0.md')
$ws.Hyperlinks.Add($ws.Range('D4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c32d239dcec73e8605290148acaa033669ff389f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf', "", "", '2f4bd7cb-c696-4503-99c3-da8190232af0.b64a3c588d434ab4c7a8a27f5428149657c67613.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/54a8ca9bebda7d35e8d6e69eac6b9ed3df74a6c7/e2e/34ded686-4006-40a0-a24a-57ef94237596.md', "", "", '2f4bd7cb-c696-4503-99c3-da8190232af0.md')
$ws.Hyperlinks.Add($ws.Range('G4'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/8d1000a3019e2fd70ed23b4162a640027ce1539d/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/34ded686-4006-40a0-a24a-57ef94237596.5c6b38f4de7feb44bf2d40b2852ecd473571ec83.de-de.xlf', "", "", '2f4bd7cb-c696-4503-99c3-da8190232af0.b64a3c588d434ab4c7a8a27f5428149657c67613.de-de.xlf')

